# Applies the "Horarios actualizados Linea 141 - 276" refresh to all three sheets.
# The source data was re-scraped (Ultima actualizacion: 13:53:08 -> 14:19:48),
# which updated several already-scraped rows (Hora_Scrap / Linea / Minutos) and
# appended newly observed arrivals at the bottom of each sheet.

$wb = $excel.ActiveWorkbook

# ===================== Sheet: LP1912 =====================
$ws = $wb.Worksheets.Item('LP1912')

# Header / summary cells
$ws.Range("A2").Value = 'Última actualización: 14:19:48'
$ws.Range("A3").Value = 'Total filas: 296'

# Data rows (modified + newly appended)
# Row 42
$ws.Range("A42").Value = '06:57:30'
$ws.Range("C42").Value = '14_ABASTO'
$ws.Range("D42").Value = 1
# Row 43
$ws.Range("A43").Value = '05:18:56'
$ws.Range("C43").Value = '10_OLMOS'
$ws.Range("D43").Value = 100
# Row 96
$ws.Range("A96").Value = '08:21:50'
$ws.Range("C96").Value = '23_HERNANDEZ'
$ws.Range("D96").Value = 40
# Row 97
$ws.Range("A97").Value = '07:20:40'
$ws.Range("C97").Value = '215A_EL PATO'
$ws.Range("D97").Value = 101
# Row 189
$ws.Range("A189").Value = '12:01:50'
$ws.Range("C189").Value = '16_SANTA ANA'
$ws.Range("D189").Value = 5
# Row 190
$ws.Range("A190").Value = '12:01:50'
$ws.Range("C190").Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Range("D190").Value = 5
# Row 191
$ws.Range("A191").Value = '10:56:30'
$ws.Range("C191").Value = '14_ABASTO'
$ws.Range("D191").Value = 70
# Row 192
$ws.Range("A192").Value = '10:26:41'
$ws.Range("C192").Value = '16_P MOR-SANTA ANA'
$ws.Range("D192").Value = 100
# Row 199
$ws.Range("A199").Value = '10:26:41'
$ws.Range("C199").Value = '215A_EL PATO'
$ws.Range("D199").Value = 114
# Row 200
$ws.Range("A200").Value = '10:56:30'
$ws.Range("C200").Value = '14_ABASTO'
$ws.Range("D200").Value = 84
# Row 209
$ws.Range("A209").Value = '12:37:14'
$ws.Range("C209").Value = '23_HERNANDEZ'
$ws.Range("D209").Value = 0
# Row 210
$ws.Range("A210").Value = '11:20:07'
$ws.Range("C210").Value = '27_EL RETIRO'
$ws.Range("D210").Value = 77
# Row 245
$ws.Range("A245").Value = '12:01:50'
$ws.Range("C245").Value = '215A_EL PATO'
$ws.Range("D245").Value = 109
# Row 246
$ws.Range("A246").Value = '13:19:56'
$ws.Range("C246").Value = '11_ETCHEVERRY'
$ws.Range("D246").Value = 31
# Row 253
$ws.Range("A253").Value = '13:53:08'
$ws.Range("C253").Value = '16_SANTA ANA'
$ws.Range("D253").Value = 4
# Row 254
$ws.Range("A254").Value = '12:37:14'
$ws.Range("C254").Value = '16_P MOR-167 Y 521'
$ws.Range("D254").Value = 80
# Row 255
$ws.Range("A255").Value = '12:37:14'
$ws.Range("C255").Value = '17_ROMERO'
$ws.Range("D255").Value = 87
# Row 256
$ws.Range("A256").Value = '13:19:56'
$ws.Range("C256").Value = '23_HERNANDEZ'
$ws.Range("D256").Value = 45
# Row 264
$ws.Range("A264").Value = '13:19:56'
$ws.Range("C264").Value = '26_HERNANDEZ'
$ws.Range("D264").Value = 61
# Row 265
$ws.Range("A265").Value = '12:37:14'
$ws.Range("C265").Value = '215C_EL PATO'
$ws.Range("D265").Value = 103
# Row 269
$ws.Range("A269").Value = '14:19:48'
$ws.Range("B269").Value = '14:34'
$ws.Range("C269").Value = '23_HERNANDEZ'
$ws.Range("D269").Value = 15
# Row 270
$ws.Range("A270").Value = '14:19:48'
$ws.Range("B270").Value = '14:44'
$ws.Range("C270").Value = '15_ABASTO'
$ws.Range("D270").Value = 25
# Row 271
$ws.Range("A271").Value = '14:19:48'
$ws.Range("B271").Value = '14:44'
$ws.Range("D271").Value = 25
# Row 272
$ws.Range("A272").Value = '13:53:08'
$ws.Range("B272").Value = '14:45'
$ws.Range("C272").Value = '14_ABASTO'
$ws.Range("D272").Value = 52
# Row 273
$ws.Range("A273").Value = '14:19:48'
$ws.Range("B273").Value = '14:46'
$ws.Range("C273").Value = '16_SANTA ANA'
$ws.Range("D273").Value = 27
# Row 274
$ws.Range("B274").Value = '14:49'
$ws.Range("C274").Value = '14_ABASTO'
$ws.Range("D274").Value = 90
# Row 275
$ws.Range("A275").Value = '12:55:01'
$ws.Range("B275").Value = '14:50'
$ws.Range("C275").Value = '14_ABASTO'
$ws.Range("D275").Value = 115
# Row 276
$ws.Range("B276").Value = '14:56'
$ws.Range("C276").Value = '16_P MOR-SANTA ANA'
$ws.Range("D276").Value = 97
# Row 277
$ws.Range("B277").Value = '14:57'
$ws.Range("C277").Value = '16_P MOR-SANTA ANA'
$ws.Range("D277").Value = 64
# Row 278
$ws.Range("A278").Value = '13:19:56'
$ws.Range("B278").Value = '14:58'
$ws.Range("C278").Value = '215B_EL PATO'
$ws.Range("D278").Value = 99
# Row 279
$ws.Range("B279").Value = '15:00'
$ws.Range("C279").Value = '81_EL PELIGRO'
$ws.Range("D279").Value = 101
# Row 280
$ws.Range("A280").Value = '13:19:56'
$ws.Range("B280").Value = '15:04'
$ws.Range("C280").Value = '10_OLMOS'
$ws.Range("D280").Value = 105
# Row 281
$ws.Range("B281").Value = '15:05'
$ws.Range("C281").Value = '10_OLMOS'
$ws.Range("D281").Value = 72
# Row 282
$ws.Range("A282").Value = '14:19:48'
$ws.Range("B282").Value = '15:06'
$ws.Range("C282").Value = '16_SANTA ANA'
$ws.Range("D282").Value = 47
# Row 283
$ws.Range("B283").Value = '15:10'
$ws.Range("C283").Value = '17_ROMERO'
$ws.Range("D283").Value = 77
# Row 284
$ws.Range("A284").Value = '13:19:56'
$ws.Range("B284").Value = '15:13'
$ws.Range("C284").Value = '11_ETCHEVERRY'
$ws.Range("D284").Value = 114
# Row 285
$ws.Range("B285").Value = '15:14'
$ws.Range("C285").Value = '11_ETCHEVERRY'
$ws.Range("D285").Value = 81
# Row 286
$ws.Range("A286").Value = '14:19:48'
$ws.Range("B286").Value = '15:21'
$ws.Range("C286").Value = '26_HERNANDEZ'
$ws.Range("D286").Value = 62
# Row 287
$ws.Range("B287").Value = '15:28'
$ws.Range("C287").Value = '26_HERNANDEZ'
$ws.Range("D287").Value = 95
# Row 288
$ws.Range("A288").Value = '14:19:48'
$ws.Range("B288").Value = '15:29'
$ws.Range("C288").Value = '14_ABASTO'
$ws.Range("D288").Value = 70
$ws.Range("E288").Value = 'LP1912'
# Row 289
$ws.Range("A289").Value = '13:53:08'
$ws.Range("B289").Value = '15:32'
$ws.Range("C289").Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Range("D289").Value = 99
$ws.Range("E289").Value = 'LP1912'
# Row 290
$ws.Range("A290").Value = '13:53:08'
$ws.Range("B290").Value = '15:35'
$ws.Range("C290").Value = '23_HERNANDEZ'
$ws.Range("D290").Value = 102
$ws.Range("E290").Value = 'LP1912'
# Row 291
$ws.Range("A291").Value = '13:53:08'
$ws.Range("B291").Value = '15:37'
$ws.Range("C291").Value = '10_OLMOS'
$ws.Range("D291").Value = 104
$ws.Range("E291").Value = 'LP1912'
# Row 292
$ws.Range("A292").Value = '14:19:48'
$ws.Range("B292").Value = '15:38'
$ws.Range("C292").Value = '215A_EL PATO'
$ws.Range("D292").Value = 79
$ws.Range("E292").Value = 'LP1912'
# Row 293
$ws.Range("A293").Value = '14:19:48'
$ws.Range("B293").Value = '15:38'
$ws.Range("C293").Value = '23_HERNANDEZ'
$ws.Range("D293").Value = 79
$ws.Range("E293").Value = 'LP1912'
# Row 294
$ws.Range("A294").Value = '13:53:08'
$ws.Range("B294").Value = '15:39'
$ws.Range("C294").Value = '215A_EL PATO'
$ws.Range("D294").Value = 106
$ws.Range("E294").Value = 'LP1912'
# Row 295
$ws.Range("A295").Value = '13:53:08'
$ws.Range("B295").Value = '15:44'
$ws.Range("C295").Value = '14_ABASTO'
$ws.Range("D295").Value = 111
$ws.Range("E295").Value = 'LP1912'
# Row 296
$ws.Range("A296").Value = '14:19:48'
$ws.Range("B296").Value = '15:46'
$ws.Range("C296").Value = '16_P MOR-167 Y 521'
$ws.Range("D296").Value = 87
$ws.Range("E296").Value = 'LP1912'
# Row 297
$ws.Range("A297").Value = '13:53:08'
$ws.Range("B297").Value = '15:47'
$ws.Range("C297").Value = '16_P MOR-167 Y 521'
$ws.Range("D297").Value = 114
$ws.Range("E297").Value = 'LP1912'
# Row 298
$ws.Range("A298").Value = '14:19:48'
$ws.Range("B298").Value = '15:53'
$ws.Range("C298").Value = '11_ETCHEVERRY'
$ws.Range("D298").Value = 94
$ws.Range("E298").Value = 'LP1912'
# Row 299
$ws.Range("A299").Value = '14:19:48'
$ws.Range("B299").Value = '15:55'
$ws.Range("C299").Value = '17_ROMERO'
$ws.Range("D299").Value = 96
$ws.Range("E299").Value = 'LP1912'
# Row 300
$ws.Range("A300").Value = '14:19:48'
$ws.Range("B300").Value = '15:56'
$ws.Range("C300").Value = '27_EL RETIRO'
$ws.Range("D300").Value = 97
$ws.Range("E300").Value = 'LP1912'
# Row 301
$ws.Range("A301").Value = '14:19:48'
$ws.Range("B301").Value = '16:14'
$ws.Range("C301").Value = '225_C ROCA-H SUR'
$ws.Range("D301").Value = 115
$ws.Range("E301").Value = 'LP1912'

# ===================== Sheet: LP1912-215 =====================
$ws = $wb.Worksheets.Item('LP1912-215')

# Header / summary cells
$ws.Range("A2").Value = 'Última actualización: 14:19:48'
$ws.Range("A3").Value = 'Total filas: 33'

# Data rows (modified + newly appended)
# Row 37
$ws.Range("A37").Value = '14:19:48'
$ws.Range("B37").Value = '15:38'
$ws.Range("D37").Value = 79
# Row 38
$ws.Range("A38").Value = '13:53:08'
$ws.Range("B38").Value = '15:39'
$ws.Range("C38").Value = '215A_EL PATO'
$ws.Range("D38").Value = 106
$ws.Range("E38").Value = 'LP1912'

# ===================== Sheet: 6203-6173 =====================
$ws = $wb.Worksheets.Item('6203-6173')

# Header / summary cells
$ws.Range("A2").Value = 'Última actualización: 14:19:48'
$ws.Range("A3").Value = 'Total filas: 44'

# Data rows (modified + newly appended)
# Row 49
$ws.Range("A49").Value = '14:19:48'
$ws.Range("B49").Value = '16:13'
$ws.Range("C49").Value = '215C_LA PLATA'
$ws.Range("D49").Value = 114
$ws.Range("E49").Value = 'L6203'
